$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 154, pushing the existing data
# (previously rows 154-200) down to rows 156-202.
$ws.Rows.Item(154).Insert()
$ws.Rows.Item(154).Insert()

# Populate the first new row (154) - Americana (o) / Primera entry dated 2021-12-29
$ws.Cells.Item(154, 1).Value = 2
$ws.Cells.Item(154, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(154, 3).Value = "Coquimbo"
$ws.Cells.Item(154, 4).Value = "2021-12-29"
$ws.Cells.Item(154, 5).Value = 4
$ws.Cells.Item(154, 6).Value = 100112021
$ws.Cells.Item(154, 7).Value = "Ají"
$ws.Cells.Item(154, 8).Value = "Americana (o)"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 400
$ws.Cells.Item(154, 11).Value = 16000
$ws.Cells.Item(154, 12).Value = 18000
$ws.Cells.Item(154, 13).Value = 17000
$ws.Cells.Item(154, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(154, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(154, 16).Value = 680
$ws.Cells.Item(154, 17).Value = 25
$ws.Cells.Item(154, 18).Value = "Hortaliza"

# Populate the second new row (155) - Inferno / Primera entry dated 2021-12-29
$ws.Cells.Item(155, 1).Value = 2
$ws.Cells.Item(155, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = "2021-12-29"
$ws.Cells.Item(155, 5).Value = 4
$ws.Cells.Item(155, 6).Value = 100112021
$ws.Cells.Item(155, 7).Value = "Ají"
$ws.Cells.Item(155, 8).Value = "Inferno"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 500
$ws.Cells.Item(155, 11).Value = 22000
$ws.Cells.Item(155, 12).Value = 24000
$ws.Cells.Item(155, 13).Value = 23000
$ws.Cells.Item(155, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(155, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(155, 16).Value = 920
$ws.Cells.Item(155, 17).Value = 25
$ws.Cells.Item(155, 18).Value = "Hortaliza"
